$wb = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item("metadata (2)")
$src.Copy($wb.Worksheets.Item(1))
$newSheet = $wb.Worksheets.Item(1)
$newSheet.Name = "temp"

$newSheet.Range("C1:D14").EntireColumn.Delete()
$newSheet.Range("A14:B14").EntireRow.Delete()

$newSheet.Range("B3").Value = "Pentius and jamaican banana "
$newSheet.Range("B6").Value = "https://na.op.gg/multi/query=jamicanbanana%2Cpentius"
$newSheet.Range("B11").Value = "Follow Pentius and jamaican banana for educational and high elo content on support and Ivern!"
$newSheet.Range("B10").Value = "Pentius and jamaican banana have both reached grandmaster and have thousands of Ivern games. They both also play support, but have continuously played Ivern since his release."
$newSheet.Range("B5").Value = "http://www.multitwitch.tv/jamicanbanana/pentiuslol"

$hl4 = $null
foreach ($hl in $newSheet.Hyperlinks) {
    if ($hl.Range.Row -eq 4) { $hl4 = $hl; break }
}
if ($hl4) { $hl4.Delete() }

$hl7 = $null
foreach ($hl in $newSheet.Hyperlinks) {
    if ($hl.Range.Row -eq 7) { $hl7 = $hl; break }
}
if ($hl7) { $hl7.Delete() }

$newSheet.Range("B4").Clear()
$newSheet.Range("B7").Clear()

$newSheet.Columns("B:B").AutoFit()

foreach ($ws in $wb.Worksheets) {
    Write-Host $ws.Name
}
